$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update data values
$ws.Range("C2").Value = 861
$ws.Range("D2").Value = 6240

$ws.Range("B3").Value = 6643
$ws.Range("C3").Value = 6643
$ws.Range("D3").Value = 6643

$ws.Range("B4").Value = 6643
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 6643

$ws.Range("B6").Value = 6643
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 3

$ws.Range("B7").Value = 1.465
$ws.Range("C7").Value = 0.021
$ws.Range("D7").Value = 0.019

# Update the selected cell/range to match the new active selection
$ws.Activate()
$ws.Range("D7").Select()
